# Update cryptocurrency price/volume figures in the "cryptos" sheet
# to reflect the latest GitHub Actions scrape (Tue Dec 26 07:54:06 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.773.79"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "2.239.37"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "`'113.50"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "`'268.87"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "`'0.626"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "`'0.605"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "`'46.18"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "`'0.0929"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "`'9.10"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "`'0.872"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "2.575.37"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "2.238.54"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "42.851.28"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "`'72.00"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "`'2.34"
$ws.Range("E22").Value = "  -5.87%  "
$ws.Range("D23").Value = "`'230.79"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").Value = "`'2.92"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").Value = "`'9.28"
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("D26").Value = "`'12.21"
$ws.Range("E26").Value = "  +6.80%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "`'40.26"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "`'2.24"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "`'3.28"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "`'173.69"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "`'0.0902"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "`'5.56"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "`'4.28"
$ws.Range("E35").Value = "  +8.90%  "
$ws.Range("D36").Value = "`'0.127"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "`'4.73"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").Value = "`'0.0375"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  -7.40%  "
$ws.Range("D42").Value = "`'13.14"
$ws.Range("E42").Value = "  -7.65%  "
$ws.Range("D43").Value = "`'0.231"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").Value = "`'0.999"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "`'5.65"
$ws.Range("E45").Value = "  -8.79%  "
$ws.Range("E46").Value = "  -3.48%  "
$ws.Range("D47").Value = "`'8.44"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "`'0.0988"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "`'100.39"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("E51").Value = "  +7.28%  "
